$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the top; everything shifts down by 2
# (old row 1 -> row 3, old row 2 -> row 4, etc.)
$ws.Rows("1:2").Insert()

# The inserted rows 1:2 picked up the header formatting (bold/border/centered)
# that used to live on row 1 and has now shifted to row 3. Copy that same
# formatting onto the new row 1, then strip the bold header formatting back
# off of row 3 (it should look like a normal data row again).
$ws.Range("A3:L3").Copy()
$ws.Range("A1:L1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A3:L3").ClearFormats()

# New row 1: numeric column indices 0..11
$headerIndex = 0
for ($col = 1; $col -le 12; $col++) {
    $ws.Cells.Item(1, $col).Value = $headerIndex
    $headerIndex = $headerIndex + 1
}

# New row 2: only B2 holds "Head", rest stay blank
$ws.Cells.Item(2, 2).Value = "Head"
